# Apply the UAE league 2023-2024 update:
# 1) Several adjacent match rows had their match-data columns (F:V) swapped
#    (the index/date columns A:E stay put).
# 2) Seven brand new match rows (79-85) are appended at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Swap columns F..V (6..22) between each pair of rows listed below.
# ---------------------------------------------------------------------
$swapPairs = @(11,12, 18,19, 39,40, 55,56, 58,59, 63,64, 74,75)

for ($i = 0; $i -lt $swapPairs.Count; $i += 2) {
    $rowA = $swapPairs[$i]
    $rowB = $swapPairs[$i + 1]

    for ($col = 6; $col -le 22; $col++) {
        $cellA = $ws.Cells.Item($rowA, $col)
        $cellB = $ws.Cells.Item($rowB, $col)
        $valA = $cellA.Value2
        $valB = $cellB.Value2
        $cellA.Value2 = $valB
        $cellB.Value2 = $valA
    }
}

# ---------------------------------------------------------------------
# 2) Append the seven new rows (79-85) at the bottom of the sheet.
# ---------------------------------------------------------------------
$newRows = @(
    @{ Row=79; Idx=78; Date=45283.57291666666;  F="Al Nasr";              G=2; H="Emirates Club";       I=1; J=1.41; K="18/12/2023 17:12"; L=1.33; M="23/12/2023 12:42"; N=4.92; O="18/12/2023 17:12"; P=6.09; Q="23/12/2023 13:43"; R=5.87; S="18/12/2023 17:12"; T=8.890000000000001; U="23/12/2023 13:43"; V="https://www.betexplorer.com/football/united-arab-emirates/uae-league/al-nasr-emirates-club/4znjvVi6/" },
    @{ Row=80; Idx=79; Date=45283.57291666666;  F="Al Jazira";            G=3; H="Hatta";                I=1; J=1.2;  K="18/12/2023 17:12"; L=1.21; M="23/12/2023 13:38"; N=6.84; O="18/12/2023 17:12"; P=7.18; Q="23/12/2023 13:43"; R=8.44; S="18/12/2023 17:12"; T=11.03;              U="23/12/2023 13:43"; V="https://www.betexplorer.com/football/united-arab-emirates/uae-league/al-jazira-hatta/r327ngpQ/" },
    @{ Row=81; Idx=80; Date=45283.6875;          F="Bani Yas";             G=2; H="Ittihad Kalba";        I=1; J=3.22; K="18/12/2023 17:12"; L=3.43; M="23/12/2023 16:29"; N=4.14; O="18/12/2023 17:12"; P=4;    Q="23/12/2023 16:29"; R=1.92; S="18/12/2023 17:12"; T=1.98;               U="23/12/2023 16:29"; V="https://www.betexplorer.com/football/united-arab-emirates/uae-league/bani-yas-ittihad-kalba/EogwsX7m/" },
    @{ Row=82; Idx=81; Date=45284.57291666666;  F="Al Bataeh";            G=1; H="Al Ain";               I=3; J=7.39; K="20/12/2023 14:12"; L=8.279999999999999; M="24/12/2023 13:44"; N=5.96; O="20/12/2023 14:12"; P=6.2;  Q="24/12/2023 13:44"; R=1.26; S="20/12/2023 14:12"; T=1.29;               U="24/12/2023 13:44"; V="https://www.betexplorer.com/football/united-arab-emirates/uae-league/al-bataeh-al-ain/rRpfwk7C/" },
    @{ Row=83; Idx=82; Date=45284.57291666666;  F="Al Sharjah";           G=1; H="Ajman";                I=1; J=1.4;  K="20/12/2023 16:42"; L=1.53; M="24/12/2023 13:44"; N=5.01; O="20/12/2023 16:42"; P=4.69; Q="24/12/2023 13:44"; R=5.84; S="20/12/2023 16:42"; T=5.42;               U="24/12/2023 13:44"; V="https://www.betexplorer.com/football/united-arab-emirates/uae-league/al-sharjah-ajman/ADnnuBx0/" },
    @{ Row=84; Idx=83; Date=45284.6875;          F="Khorfakkan";           G=3; H="Al Wahda";             I=2; J=6.29; K="20/12/2023 14:12"; L=6.6;  M="24/12/2023 16:25"; N=5.15; O="20/12/2023 14:12"; P=5.03; Q="24/12/2023 16:28"; R=1.36; S="20/12/2023 14:12"; T=1.43;               U="24/12/2023 16:21"; V="https://www.betexplorer.com/football/united-arab-emirates/uae-league/khorfakkan-al-wahda/zefZsDhs/" },
    @{ Row=85; Idx=84; Date=45284.6875;          F="Shabab Al-Ahli Dubai"; G=1; H="Al Wasl";              I=2; J=2.76; K="20/12/2023 16:42"; L=2.57; M="24/12/2023 16:28"; N=3.84; O="20/12/2023 16:42"; P=3.81; Q="24/12/2023 16:28"; R=2.22; S="20/12/2023 16:42"; T=2.55;               U="24/12/2023 16:28"; V="https://www.betexplorer.com/football/united-arab-emirates/uae-league/shabab-al-ahli-dubai-al-wasl/dArrtiNg/" }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value2  = $r.Idx
    $ws.Cells.Item($row, 2).Value2  = "united-arab-emirates"
    $ws.Cells.Item($row, 3).Value2  = "uae-league"
    $ws.Cells.Item($row, 4).Value2  = "2023-2024"
    $ws.Cells.Item($row, 5).Value2  = $r.Date
    $ws.Cells.Item($row, 6).Value2  = $r.F
    $ws.Cells.Item($row, 7).Value2  = $r.G
    $ws.Cells.Item($row, 8).Value2  = $r.H
    $ws.Cells.Item($row, 9).Value2  = $r.I
    $ws.Cells.Item($row, 10).Value2 = $r.J
    $ws.Cells.Item($row, 11).Value2 = $r.K
    $ws.Cells.Item($row, 12).Value2 = $r.L
    $ws.Cells.Item($row, 13).Value2 = $r.M
    $ws.Cells.Item($row, 14).Value2 = $r.N
    $ws.Cells.Item($row, 15).Value2 = $r.O
    $ws.Cells.Item($row, 16).Value2 = $r.P
    $ws.Cells.Item($row, 17).Value2 = $r.Q
    $ws.Cells.Item($row, 18).Value2 = $r.R
    $ws.Cells.Item($row, 19).Value2 = $r.S
    $ws.Cells.Item($row, 20).Value2 = $r.T
    $ws.Cells.Item($row, 21).Value2 = $r.U
    $ws.Cells.Item($row, 22).Value2 = $r.V

    # Match the formatting used by the existing rows: bold/centered/bordered
    # index column (A) and date-time formatted column (E).
    $ws.Range("A78").Copy()
    $ws.Range("A" + $row).PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    $ws.Range("E78").Copy()
    $ws.Range("E" + $row).PasteSpecial(-4122)
    $excel.CutCopyMode = $false
}
